# "uniform the load battlefile mechanical"
# Reorder the semicolon separated tokens in the BattleMap data rows so the
# map id (510180010 / 51018001) is always the leading token, matching the
# new uniform parsing order used by the battle-file loader.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Id=1)
$ws.Range("D4").Value = "510180010;4;;51018002;2;1;51018002;2;7"
$ws.Range("E4").Value = "51018001;19;4;51018002;17;1;51018002;17;7"

# Row 5 (Id=2)
$ws.Range("D5").Value = "51018001;1;2"

# Row 6 (Id=3)
$ws.Range("D6").Value = "51018001;0;2"
$ws.Range("E6").Value = "51018001;10;2"

# back to row 5
$ws.Range("E5").Value = "51018001;9;2"

# Row 7 (Id=4)
$ws.Range("D7").Value = "51018001;8;4"
$ws.Range("E7").Value = "51018001;14;4"

# Widen columns D and E so the longer strings are fully visible.
# (ColumnWidth is expressed in "characters"; the runtime quantizes it to
# whole pixels on save, so we pick the values whose resulting OOXML
# <col width> is exactly/closest to the target 49.25 / 58.)
$ws.Columns.Item(4).ColumnWidth = 48.57142857142857
$ws.Columns.Item(5).ColumnWidth = 57.285714285714285

# Move the active selection to E7, the last cell touched.
$ws.Range("E7").Select()
